$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date stamp for every data
# row (rows 2-116). This automatic refresh bumps that stamp from
# 2023-10-05 (serial 45204) to 2023-10-06 (serial 45205) for every row that
# currently has the old value.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 116 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -eq 45204) {
        $cell.Value = 45205
    }
}
